# Apply the vocabulary.xlsx edit:
#  - rename the "zonmwpc" vocabulary/prefix to "covid19" (new PURL base URI)
#  - replace every "zonmwpc:" identifier reference with "covid19:" throughout the sheet
#  - overwrite rows 260-262 (former "prevention/aftercare/free-from-disease phase" rows)
#    with the new iadopt "variable"/"property"/"constraint" concepts
#  - retarget row 263 from "fungus" to "genomics data" (and update its broader concept)
#  - retarget the "skos:broader" values of rows 264-265 ("perineal/cloacal swab")
#  - drop the three rows that used to hold the "variable"/"property"/"constraint"
#    concepts at the bottom of the sheet (266-268), since that content moved up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) ConceptScheme URI (row 1) and PREFIX declaration (row 2)
$ws.Range("B1").Value = "http://purl.org/zonmw/covid19/"
$ws.Range("B2").Value = "covid19"
$ws.Range("C2").Value = "http://purl.org/zonmw/covid19/"

# 2) Globally rename every "zonmwpc:xxxxx" identifier/reference to "covid19:xxxxx"
$used = $ws.UsedRange
$null = $used.Replace("zonmwpc:", "covid19:")

# 3) Row 260: "prevention phase" -> "variable"
$ws.Range("A260").Value = "covid19:10241"
$ws.Range("B260").Value = "variable"
$ws.Range("C260").Value = ""
$ws.Range("D260").Value = ""
$ws.Range("E260").Value = "A description of something observed or derived, minimally consisting of  an ObjectOfInterest and its Property."
$ws.Range("F260").Value = ""
$ws.Range("G260").Value = "https://w3id.org/iadopt/ont/iadopt/variable"
$ws.Range("H260").Value = ""
$ws.Range("I260").Value = ""
$ws.Range("J260").Value = "https://w3id.org/iadopt/ont/iadopt/variable"

# 4) Row 261: "aftercare phase" -> "property"
$ws.Range("A261").Value = "covid19:10242"
$ws.Range("B261").Value = "property"
$ws.Range("C261").Value = ""
$ws.Range("D261").Value = ""
$ws.Range("E261").Value = "A type of a characteristic of the ObjectOfInterest."
$ws.Range("F261").Value = ""
$ws.Range("G261").Value = "https://w3id.org/iadopt/ont/iadopt/property"
$ws.Range("H261").Value = ""
$ws.Range("I261").Value = ""
$ws.Range("J261").Value = "https://w3id.org/iadopt/ont/iadopt/property"

# 5) Row 262: "free-from-disease phase" -> "constraint"
$ws.Range("A262").Value = "covid19:10243"
$ws.Range("B262").Value = "constraint"
$ws.Range("C262").Value = ""
$ws.Range("D262").Value = ""
$ws.Range("E262").Value = "A Constraint limits the scope of the observation and confines the context to a particular state. It describes properties of the involved entities that are relevant to the particular observation."
$ws.Range("F262").Value = ""
$ws.Range("G262").Value = "https://w3id.org/iadopt/ont/iadopt/constraint"
$ws.Range("H262").Value = ""
$ws.Range("I262").Value = ""
$ws.Range("J262").Value = "https://w3id.org/iadopt/ont/iadopt/constraint"

# 6) Row 263: "fungus" -> "genomics data" / "sequencing data", new broader concept
$ws.Range("B263").Value = "genomics data"
$ws.Range("C263").Value = "sequencing data"
$ws.Range("H263").Value = "covid19:10107"

# 7) Rows 264-265 ("perineal swab" / "cloacal swab") get a new broader concept
$ws.Range("H264").Value = "covid19:10144"
$ws.Range("H265").Value = "covid19:10144"

# 8) Remove the now-duplicated rows 266-268 that used to hold the
#    variable/property/constraint concepts (their content was moved to 260-262)
$ws.Range("A266:A268").EntireRow.Delete()
